# Add a new "2022-Q3" sheet with fund-holding detail data, insert it right
# after "总计" (becoming the 2nd sheet, before "2022-Q2"), and update the
# "总计" summary sheet with the new quarter's totals (shifting the other
# rows down by one).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Update the "总计" (totals) summary sheet: insert a new row for the
#    2022-Q3 quarter at the top of the data (row 2), push everything else
#    down by one row.
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")
$summary.Rows.Item(2).Insert()

$summary.Cells.Item(2, 1).Value = 0
$summary.Cells.Item(2, 2).Value = "2022-Q3"
$summary.Cells.Item(2, 3).Value = 8
$summary.Cells.Item(2, 4).Value = 1.42

# Renumber the "序号" (index) column A for the rows that got pushed down.
$summary.Cells.Item(3, 1).Value = 1
$summary.Cells.Item(4, 1).Value = 2
$summary.Cells.Item(5, 1).Value = 3
$summary.Cells.Item(6, 1).Value = 4

# ---------------------------------------------------------------------
# 2) Insert the new "2022-Q3" detail worksheet right before "2022-Q2"
#    (i.e. as the 2nd tab, right after "总计").
# ---------------------------------------------------------------------
$beforeSheet = $wb.Worksheets.Item("2022-Q2")
$q3 = $wb.Worksheets.Add($beforeSheet)
$q3.Name = "2022-Q3"

$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $cell = $q3.Cells.Item(1, $i + 2)
    $cell.Value = $headers[$i]
}
$headerRange = $q3.Range("B1:H1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

# Columns B..G hold text (fund code / name / scale / position / ratio /
# market value) in the source data, so force Text format before writing
# so Excel doesn't "helpfully" reinterpret them as numbers (which would
# strip leading/trailing zeros such as "012860" or "0.0680").
$q3.Range("B2:G9").NumberFormat = "@"

$rows = @(
    @("513500", "博时标普500ETF（QDII）",              "71.37", "96.44", "1.54", "1.0991", 7),
    @("161125", "易方达标普500指数（QDII-LOF）人民币",   "4.74", "90.72", "1.46", "0.0692", 7),
    @("012860", "易方达标普500指数（QDII-LOF）人民币 C", "4.74", "90.72", "1.46", "0.0692", 7),
    @("003718", "易方达标普500指数（QDII-LOF）美元A",    "4.66", "90.72", "1.46", "0.0680", 7),
    @("013329", "嘉实全球价值股票（QDII）美元现汇",       "1.68", "90.63", "3.22", "0.0541", 3),
    @("013328", "嘉实全球价值股票（QDII）人民币",         "1.68", "90.63", "3.22", "0.0541", 3),
    @("159612", "国泰标普500ETF（QDII）",                "0.55", "91.40", "1.45", "0.0080", 7),
    @("012861", "易方达标普500指数（QDII-LOF）美元 C",    "0.08", "90.72", "1.46", "0.0012", 7)
)

for ($r = 0; $r -lt $rows.Length; $r++) {
    $row = $rows[$r]
    $excelRow = $r + 2
    $q3.Cells.Item($excelRow, 1).Value = $r
    $q3.Cells.Item($excelRow, 2).Value = $row[0]
    $q3.Cells.Item($excelRow, 3).Value = $row[1]
    $q3.Cells.Item($excelRow, 4).Value = $row[2]
    $q3.Cells.Item($excelRow, 5).Value = $row[3]
    $q3.Cells.Item($excelRow, 6).Value = $row[4]
    $q3.Cells.Item($excelRow, 7).Value = $row[5]
    $q3.Cells.Item($excelRow, 8).Value = $row[6]
}

$q3.Range("A2:A9").Font.Bold = $true
